$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 'ba'
$ws.Range("J3").Value = 'Appreciation'
$ws.Range("I5").Value = 'sv'
$ws.Range("J5").Value = 'Statement-opinion'
$ws.Range("I8").Value = 'sd'
$ws.Range("J8").Value = 'Statement-non-opinion'
$ws.Range("I16").Value = 'qy'
$ws.Range("J16").Value = 'Yes-No-Question'
$ws.Range("I23").Value = 'b'
$ws.Range("J23").Value = 'Acknowledge (Backchannel)'
$ws.Range("I27").Value = 'sd'
$ws.Range("J27").Value = 'Statement-non-opinion'
$ws.Range("I33").Value = 'sd'
$ws.Range("J33").Value = 'Statement-non-opinion'
$ws.Range("I35").Value = 'sd'
$ws.Range("J35").Value = 'Statement-non-opinion'
$ws.Range("I36").Value = 'sd'
$ws.Range("J36").Value = 'Statement-non-opinion'
$ws.Range("I40").Value = 'sd'
$ws.Range("J40").Value = 'Statement-non-opinion'
$ws.Range("I44").Value = 'b'
$ws.Range("J44").Value = 'Acknowledge (Backchannel)'
$ws.Range("I45").Value = 'b'
$ws.Range("J45").Value = 'Acknowledge (Backchannel)'
$ws.Range("I47").Value = 'sd'
$ws.Range("J47").Value = 'Statement-non-opinion'
$ws.Range("I51").Value = '%'
$ws.Range("J51").Value = 'Uninterpretable'
$ws.Range("I54").Value = 'aa'
$ws.Range("J54").Value = 'Agree/Accept'
$ws.Range("I58").Value = 'sd'
$ws.Range("J58").Value = 'Statement-non-opinion'
$ws.Range("I62").Value = 'ba'
$ws.Range("J62").Value = 'Appreciation'
$ws.Range("I63").Value = 'b'
$ws.Range("J63").Value = 'Acknowledge (Backchannel)'
$ws.Range("I64").Value = 'qy'
$ws.Range("J64").Value = 'Yes-No-Question'
$ws.Range("I65").Value = 'sd'
$ws.Range("J65").Value = 'Statement-non-opinion'
$ws.Range("I70").Value = 'sd'
$ws.Range("J70").Value = 'Statement-non-opinion'
$ws.Range("I72").Value = 'sd'
$ws.Range("J72").Value = 'Statement-non-opinion'
$ws.Range("I78").Value = 'sd'
$ws.Range("J78").Value = 'Statement-non-opinion'
$ws.Range("I80").Value = 'aa'
$ws.Range("J80").Value = 'Agree/Accept'
$ws.Range("I88").Value = 'sd'
$ws.Range("J88").Value = 'Statement-non-opinion'
$ws.Range("I92").Value = 'b'
$ws.Range("J92").Value = 'Acknowledge (Backchannel)'
$ws.Range("I97").Value = 'sd'
$ws.Range("J97").Value = 'Statement-non-opinion'
$ws.Range("I98").Value = 'sd'
$ws.Range("J98").Value = 'Statement-non-opinion'
$ws.Range("I102").Value = 'sd'
$ws.Range("J102").Value = 'Statement-non-opinion'
$ws.Range("I115").Value = 'sv'
$ws.Range("J115").Value = 'Statement-opinion'
$ws.Range("I118").Value = 'sd'
$ws.Range("J118").Value = 'Statement-non-opinion'
$ws.Range("I127").Value = 'b'
$ws.Range("J127").Value = 'Acknowledge (Backchannel)'
$ws.Range("I134").Value = 'b'
$ws.Range("J134").Value = 'Acknowledge (Backchannel)'
$ws.Range("I143").Value = 'ba'
$ws.Range("J143").Value = 'Appreciation'
$ws.Range("I149").Value = 'sv'
$ws.Range("J149").Value = 'Statement-opinion'
$ws.Range("I150").Value = 'sv'
$ws.Range("J150").Value = 'Statement-opinion'
$ws.Range("I166").Value = 'sd'
$ws.Range("J166").Value = 'Statement-non-opinion'
$ws.Range("I173").Value = 'sd'
$ws.Range("J173").Value = 'Statement-non-opinion'
$ws.Range("I185").Value = 'sd'
$ws.Range("J185").Value = 'Statement-non-opinion'
$ws.Range("I210").Value = 'b'
$ws.Range("J210").Value = 'Acknowledge (Backchannel)'
$ws.Range("I215").Value = 'sd'
$ws.Range("J215").Value = 'Statement-non-opinion'
$ws.Range("I221").Value = 'sd'
$ws.Range("J221").Value = 'Statement-non-opinion'
$ws.Range("I224").Value = 'sd'
$ws.Range("J224").Value = 'Statement-non-opinion'
$ws.Range("I228").Value = 'sd'
$ws.Range("J228").Value = 'Statement-non-opinion'
$ws.Range("I236").Value = 'sv'
$ws.Range("J236").Value = 'Statement-opinion'
$ws.Range("I274").Value = 'b'
$ws.Range("J274").Value = 'Acknowledge (Backchannel)'
$ws.Range("I277").Value = 'aa'
$ws.Range("J277").Value = 'Agree/Accept'
$ws.Range("I284").Value = 'sv'
$ws.Range("J284").Value = 'Statement-opinion'
$ws.Range("I291").Value = 'sd'
$ws.Range("J291").Value = 'Statement-non-opinion'
$ws.Range("I294").Value = 'sd'
$ws.Range("J294").Value = 'Statement-non-opinion'
$ws.Range("I295").Value = 'sv'
$ws.Range("J295").Value = 'Statement-opinion'
$ws.Range("I315").Value = 'sd'
$ws.Range("J315").Value = 'Statement-non-opinion'
$ws.Range("I319").Value = '%'
$ws.Range("J319").Value = 'Uninterpretable'
$ws.Range("I325").Value = 'sd'
$ws.Range("J325").Value = 'Statement-non-opinion'
$ws.Range("I327").Value = 'sd'
$ws.Range("J327").Value = 'Statement-non-opinion'
$ws.Range("I328").Value = 'b'
$ws.Range("J328").Value = 'Acknowledge (Backchannel)'
$ws.Range("I330").Value = 'sv'
$ws.Range("J330").Value = 'Statement-opinion'
$ws.Range("I332").Value = 'aa'
$ws.Range("J332").Value = 'Agree/Accept'
$ws.Range("I341").Value = 'sv'
$ws.Range("J341").Value = 'Statement-opinion'
$ws.Range("I347").Value = 'sv'
$ws.Range("J347").Value = 'Statement-opinion'
$ws.Range("I350").Value = 'sd'
$ws.Range("J350").Value = 'Statement-non-opinion'
$ws.Range("I351").Value = 'aa'
$ws.Range("J351").Value = 'Agree/Accept'
$ws.Range("I354").Value = 'sv'
$ws.Range("J354").Value = 'Statement-opinion'
$ws.Range("I360").Value = 'sd'
$ws.Range("J360").Value = 'Statement-non-opinion'
$ws.Range("I367").Value = 'sd'
$ws.Range("J367").Value = 'Statement-non-opinion'
$ws.Range("I380").Value = 'b'
$ws.Range("J380").Value = 'Acknowledge (Backchannel)'
$ws.Range("I404").Value = 'b'
$ws.Range("J404").Value = 'Acknowledge (Backchannel)'
$ws.Range("I407").Value = 'sd'
$ws.Range("J407").Value = 'Statement-non-opinion'
$ws.Range("I408").Value = 'sd'
$ws.Range("J408").Value = 'Statement-non-opinion'
$ws.Range("I415").Value = 'sd'
$ws.Range("J415").Value = 'Statement-non-opinion'
$ws.Range("I417").Value = 'aa'
$ws.Range("J417").Value = 'Agree/Accept'
$ws.Range("I420").Value = 'sd'
$ws.Range("J420").Value = 'Statement-non-opinion'
$ws.Range("I429").Value = 'sd'
$ws.Range("J429").Value = 'Statement-non-opinion'
$ws.Range("I432").Value = '%'
$ws.Range("J432").Value = 'Uninterpretable'
$ws.Range("I442").Value = 'sv'
$ws.Range("J442").Value = 'Statement-opinion'
$ws.Range("I449").Value = 'sv'
$ws.Range("J449").Value = 'Statement-opinion'
$ws.Range("I457").Value = 'sv'
$ws.Range("J457").Value = 'Statement-opinion'
$ws.Range("I459").Value = 'sv'
$ws.Range("J459").Value = 'Statement-opinion'
$ws.Range("I461").Value = 'sv'
$ws.Range("J461").Value = 'Statement-opinion'
$ws.Range("I463").Value = 'ba'
$ws.Range("J463").Value = 'Appreciation'
